$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "description" column (C) entirely; "facilitators" (old D) shifts left into C.
$ws.Columns("C").Delete()

# Fix header typo: "Session" -> "SessIOn"
$ws.Range("B1").Value = "SessIOn"
